$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 0
$wsExhibit.Range("F8").Value = 144
$wsExhibit.Range("F9").Value = 63

# Sheet "全部类型" (all types) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 117
$wsAll.Range("F4").Value = 0
$wsAll.Range("F5").Value = 17
$wsAll.Range("F9").Value = 63
$wsAll.Range("F10").Value = 493
